$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 53
$ws_ALC.Range("H53").Value = 125.46154
$ws_ALC.Range("I53").Value = 137.36363
$ws_ALC.Range("J53").Value = 60
$ws_ALC.Range("K53").Value = 137.36363
$ws_ALC.Range("L53").Value = 60
$ws_ALC.Range("M53").Value = 499.63637
$ws_ALC.Range("N53").Value = -1334

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws_ARM.Range("H2").Value = 1899.9
$ws_ARM.Range("I2").Value = 1049.75
$ws_ARM.Range("J2").Value = 2466.6667
$ws_ARM.Range("K2").Value = 1049.75
$ws_ARM.Range("L2").Value = 2466.6667
$ws_ARM.Range("M2").Value = -936.75
$ws_ARM.Range("N2").Value = -2692.6667

# ARM row 32
$ws_ARM.Range("H32").Value = 5473.353
$ws_ARM.Range("I32").Value = 3454.4893
$ws_ARM.Range("J32").Value = 29195
$ws_ARM.Range("K32").Value = 3454.4893
$ws_ARM.Range("L32").Value = 29195
$ws_ARM.Range("M32").Value = -3167.4893
$ws_ARM.Range("N32").Value = -29769

# ARM row 110
$ws_ARM.Range("H110").Value = 2245.6667
$ws_ARM.Range("I110").Value = 1022.2
$ws_ARM.Range("J110").Value = 3775
$ws_ARM.Range("K110").Value = 1022.2
$ws_ARM.Range("L110").Value = 3775
$ws_ARM.Range("M110").Value = 1022.8
$ws_ARM.Range("N110").Value = -7865

# ARM row 116
$ws_ARM.Range("H116").Value = 1899.9
$ws_ARM.Range("I116").Value = 1049.75
$ws_ARM.Range("J116").Value = 2466.6667
$ws_ARM.Range("K116").Value = 1049.75
$ws_ARM.Range("L116").Value = 2466.6667
$ws_ARM.Range("M116").Value = 1244.25
$ws_ARM.Range("N116").Value = -7054.6667

# ARM row 132
$ws_ARM.Range("H132").Value = 4957
$ws_ARM.Range("I132").Value = 0
$ws_ARM.Range("J132").Value = 4957
$ws_ARM.Range("K132").Value = 0
$ws_ARM.Range("L132").Value = 14871
$ws_ARM.Range("N132").Value = -19931

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws_BSM.Range("H3").Value = 1899.9
$ws_BSM.Range("I3").Value = 1049.75
$ws_BSM.Range("J3").Value = 2466.6667
$ws_BSM.Range("K3").Value = 1049.75
$ws_BSM.Range("L3").Value = 2466.6667
$ws_BSM.Range("M3").Value = -935.75
$ws_BSM.Range("N3").Value = -2694.6667

# BSM row 80
$ws_BSM.Range("H80").Value = 460.14816
$ws_BSM.Range("I80").Value = 75.14286
$ws_BSM.Range("J80").Value = 594.9
$ws_BSM.Range("K80").Value = 75.14286
$ws_BSM.Range("L80").Value = 594.9
$ws_BSM.Range("M80").Value = 922.85714
$ws_BSM.Range("N80").Value = -2590.9

# BSM row 83
$ws_BSM.Range("H83").Value = 460.14816
$ws_BSM.Range("I83").Value = 75.14286
$ws_BSM.Range("J83").Value = 594.9
$ws_BSM.Range("K83").Value = 375.7143
$ws_BSM.Range("L83").Value = 2974.5
$ws_BSM.Range("M83").Value = 4616.2857
$ws_BSM.Range("N83").Value = -12958.5

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 134
$ws_CRP.Range("H134").Value = 2248.5667
$ws_CRP.Range("I134").Value = 1848.9048
$ws_CRP.Range("J134").Value = 3181.111
$ws_CRP.Range("K134").Value = 5546.7144
$ws_CRP.Range("L134").Value = 9543.332999999999
$ws_CRP.Range("M134").Value = -3011.7144
$ws_CRP.Range("N134").Value = -14613.333

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 8
$ws_CUL.Range("H8").Value = 54.333332
$ws_CUL.Range("I8").Value = 54.333332
$ws_CUL.Range("J8").Value = 0
$ws_CUL.Range("K8").Value = 162.999996
$ws_CUL.Range("L8").Value = 0
$ws_CUL.Range("M8").Value = -23.99999600000001

# CUL row 23
$ws_CUL.Range("H23").Value = 215.86957
$ws_CUL.Range("I23").Value = 229.83333
$ws_CUL.Range("J23").Value = 210.94118
$ws_CUL.Range("K23").Value = 689.49999
$ws_CUL.Range("L23").Value = 632.82354
$ws_CUL.Range("M23").Value = -454.49999
$ws_CUL.Range("N23").Value = -1102.82354

# CUL row 35
$ws_CUL.Range("H35").Value = 1325
$ws_CUL.Range("I35").Value = 800
$ws_CUL.Range("J35").Value = 1500
$ws_CUL.Range("K35").Value = 2400
$ws_CUL.Range("L35").Value = 4500
$ws_CUL.Range("M35").Value = -2112
$ws_CUL.Range("N35").Value = -5076

# CUL row 55
$ws_CUL.Range("H55").Value = 33817.08
$ws_CUL.Range("I55").Value = 156245.4
$ws_CUL.Range("J55").Value = 3210
$ws_CUL.Range("K55").Value = 468736.2
$ws_CUL.Range("L55").Value = 9630
$ws_CUL.Range("M55").Value = -468559.2
$ws_CUL.Range("N55").Value = -9984

# CUL row 68
$ws_CUL.Range("H68").Value = 4028.8572
$ws_CUL.Range("I68").Value = 0
$ws_CUL.Range("J68").Value = 4028.8572
$ws_CUL.Range("K68").Value = 0
$ws_CUL.Range("L68").Value = 12086.5716
$ws_CUL.Range("M68").Value = $null
$ws_CUL.Range("N68").Value = -13708.5716

# CUL row 71
$ws_CUL.Range("H71").Value = 4028.8572
$ws_CUL.Range("I71").Value = 0
$ws_CUL.Range("J71").Value = 4028.8572
$ws_CUL.Range("K71").Value = 0
$ws_CUL.Range("L71").Value = 36259.7148
$ws_CUL.Range("M71").Value = $null
$ws_CUL.Range("N71").Value = -44371.7148

# CUL row 80
$ws_CUL.Range("H80").Value = 1763.125
$ws_CUL.Range("I80").Value = 801
$ws_CUL.Range("J80").Value = 2083.8333
$ws_CUL.Range("K80").Value = 2403
$ws_CUL.Range("L80").Value = 6251.499899999999
$ws_CUL.Range("M80").Value = -1467
$ws_CUL.Range("N80").Value = -8123.499899999999

# CUL row 83
$ws_CUL.Range("H83").Value = 1763.125
$ws_CUL.Range("I83").Value = 801
$ws_CUL.Range("J83").Value = 2083.8333
$ws_CUL.Range("K83").Value = 7209
$ws_CUL.Range("L83").Value = 18754.4997
$ws_CUL.Range("M83").Value = -2529
$ws_CUL.Range("N83").Value = -28114.4997

# CUL row 97
$ws_CUL.Range("H97").Value = 348.8889
$ws_CUL.Range("I97").Value = 235
$ws_CUL.Range("J97").Value = 440
$ws_CUL.Range("K97").Value = 705
$ws_CUL.Range("L97").Value = 1320
$ws_CUL.Range("M97").Value = -209
$ws_CUL.Range("N97").Value = -2312

# CUL row 109
$ws_CUL.Range("H109").Value = 98.666664
$ws_CUL.Range("I109").Value = 98.666664
$ws_CUL.Range("J109").Value = 0
$ws_CUL.Range("K109").Value = 295.999992
$ws_CUL.Range("L109").Value = 0
$ws_CUL.Range("M109").Value = 744.000008
$ws_CUL.Range("N109").Value = $null

# CUL row 113
$ws_CUL.Range("H113").Value = 402.60294
$ws_CUL.Range("I113").Value = 362.17242
$ws_CUL.Range("J113").Value = 432.66666
$ws_CUL.Range("K113").Value = 1086.51726
$ws_CUL.Range("L113").Value = 1297.99998
$ws_CUL.Range("M113").Value = 1083.48274
$ws_CUL.Range("N113").Value = -5637.999980000001

# CUL row 122
$ws_CUL.Range("H122").Value = 5205.522
$ws_CUL.Range("I122").Value = 530.8182
$ws_CUL.Range("J122").Value = 9490.666999999999
$ws_CUL.Range("K122").Value = 4777.3638
$ws_CUL.Range("L122").Value = 85416.003
$ws_CUL.Range("M122").Value = -2327.3638
$ws_CUL.Range("N122").Value = -90316.003

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 103
$ws_GSM.Range("H103").Value = 30000
$ws_GSM.Range("I103").Value = 0
$ws_GSM.Range("J103").Value = 30000
$ws_GSM.Range("K103").Value = 0
$ws_GSM.Range("L103").Value = 30000
$ws_GSM.Range("N103").Value = -32344

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws_LTW.Range("H7").Value = 1874.4193
$ws_LTW.Range("I7").Value = 1864.3572
$ws_LTW.Range("J7").Value = 1968.3334
$ws_LTW.Range("K7").Value = 1864.3572
$ws_LTW.Range("L7").Value = 1968.3334
$ws_LTW.Range("M7").Value = -1752.3572

# LTW row 22
$ws_LTW.Range("H22").Value = 381.46667
$ws_LTW.Range("I22").Value = 378.89474
$ws_LTW.Range("J22").Value = 395.42856
$ws_LTW.Range("K22").Value = 378.89474
$ws_LTW.Range("L22").Value = 395.42856
$ws_LTW.Range("M22").Value = -83.89474000000001
$ws_LTW.Range("N22").Value = -985.4285600000001

# LTW row 27
$ws_LTW.Range("H27").Value = 381.46667
$ws_LTW.Range("I27").Value = 378.89474
$ws_LTW.Range("J27").Value = 395.42856
$ws_LTW.Range("K27").Value = 378.89474
$ws_LTW.Range("L27").Value = 395.42856
$ws_LTW.Range("M27").Value = -271.89474
$ws_LTW.Range("N27").Value = -609.4285600000001

# LTW row 46
$ws_LTW.Range("H46").Value = 1367
$ws_LTW.Range("I46").Value = 1100.5
$ws_LTW.Range("J46").Value = 1900
$ws_LTW.Range("K46").Value = 1100.5
$ws_LTW.Range("L46").Value = 1900
$ws_LTW.Range("M46").Value = -912.5
$ws_LTW.Range("N46").Value = -2276

# LTW row 55
$ws_LTW.Range("H55").Value = 165.11111
$ws_LTW.Range("I55").Value = 138.03448
$ws_LTW.Range("J55").Value = 277.2857
$ws_LTW.Range("K55").Value = 138.03448
$ws_LTW.Range("L55").Value = 277.2857
$ws_LTW.Range("M55").Value = 34.96552
$ws_LTW.Range("N55").Value = -623.2857

# LTW row 126
$ws_LTW.Range("H126").Value = 1874.4193
$ws_LTW.Range("I126").Value = 1864.3572
$ws_LTW.Range("J126").Value = 1968.3334
$ws_LTW.Range("K126").Value = 5593.071599999999
$ws_LTW.Range("L126").Value = 5905.0002
$ws_LTW.Range("M126").Value = -3123.071599999999
